$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-06-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-06 Friday", 2) | Out-Null

# Update each arithmetic expression cell in the table (row-major order)
$t = $d.Tables.Item(1)
$values = @(
    "14-6=",
    "71+24=",
    "7+53=",
    "80-3=",
    "92-34=",
    "80-56=",
    "44-15=",
    "28+20=",
    "81-70=",
    "21+18=",
    "11+86=",
    "67-38=",
    "49+43=",
    "69-28=",
    "89-9=",
    "71-31=",
    "52-23=",
    "12-1=",
    "63-61=",
    "78-8=",
    "0+56=",
    "52-4=",
    "54+40=",
    "16-14=",
    "22-8=",
    "2+38=",
    "11+43=",
    "76-68=",
    "55+44=",
    "96-17=",
    "62+31=",
    "38-27=",
    "84-56=",
    "30-19=",
    "25+15=",
    "4+11=",
    "36+11=",
    "18+15=",
    "4+2=",
    "14+22=",
    "27+50=",
    "59+11=",
    "86+4=",
    "30+13=",
    "82-63=",
    "74+9=",
    "51+4=",
    "30+30=",
    "40+5=",
    "52-5=",
    "38-22=",
    "58-9=",
    "43+14=",
    "61-33=",
    "58+22=",
    "86-8=",
    "30+23=",
    "32+47=",
    "83-78=",
    "99-81=",
    "5+79=",
    "54-26=",
    "9+64=",
    "44+54=",
    "95-25=",
    "98-52=",
    "38+40=",
    "40+42=",
    "72-1=",
    "67+32=",
    "62-23=",
    "20+20=",
    "18+28=",
    "72-26=",
    "77-69=",
    "23+56=",
    "0+50=",
    "66-35=",
    "79-8=",
    "56+6=",
    "81-53=",
    "33-17=",
    "46+11=",
    "96-33=",
    "76-51=",
    "25+67=",
    "54+23=",
    "40-26=",
    "26-4=",
    "75-42=",
    "16+55=",
    "10+55=",
    "11+74=",
    "46-9=",
    "96-41=",
    "76-51=",
    "48+19=",
    "25+18=",
    "69-50=",
    "7-1="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
if (($rows * $cols) -ne $values.Count) {
    throw ("Table shape " + $rows + "x" + $cols + " does not match " + $values.Count + " replacement values")
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")
